$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pecam1"
$ws.Range("C2").Value = "Cd38"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 234.5813823333333
$ws.Range("H2").Value = 703.744147
$ws.Range("I2").Value = 0.9782746109134588
$ws.Range("J2").Value = 0.9782746109134588
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 40.61064433333333
$ws.Range("N2").Value = 121.831933
$ws.Range("O2").Value = 0.9365419382944963
$ws.Range("P2").Value = 0.9365419382944962
$ws.Range("Q2").Value = 9526.501085160682
$ws.Range("R2").Value = 85738.50976644615
$ws.Range("S2").Value = 0.9161952002891849
$ws.Range("T2").Value = 0.9161952002891848
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pecam1"
$ws.Range("C3").Value = "Cd38"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 234.5813823333333
$ws.Range("H3").Value = 703.744147
$ws.Range("I3").Value = 0.9782746109134588
$ws.Range("J3").Value = 0.9782746109134588
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.454549333333333
$ws.Range("N3").Value = 7.363648
$ws.Range("O3").Value = 0.05660556309845621
$ws.Range("P3").Value = 0.05660556309845621
$ws.Range("Q3").Value = 575.791575618695
$ws.Range("R3").Value = 5182.124180568256
$ws.Range("S3").Value = 0.05537578521567949
$ws.Range("T3").Value = 0.05537578521567948
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pecam1"
$ws.Range("C4").Value = "Cd38"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 234.5813823333333
$ws.Range("H4").Value = 703.744147
$ws.Range("I4").Value = 0.9782746109134588
$ws.Range("J4").Value = 0.9782746109134588
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2971403333333333
$ws.Range("N4").Value = 0.891421
$ws.Range("O4").Value = 0.006852498607047613
$ws.Range("P4").Value = 0.006852498607047612
$ws.Range("Q4").Value = 69.70359014032077
$ws.Range("R4").Value = 627.3323112628871
$ws.Range("S4").Value = 0.006703625408594522
$ws.Range("T4").Value = 0.006703625408594521
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pecam1"
$ws.Range("C5").Value = "Cd38"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.110333333333333
$ws.Range("H5").Value = 6.331
$ws.Range("I5").Value = 0.008800721950008783
$ws.Range("J5").Value = 0.008800721950008783
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 40.61064433333333
$ws.Range("N5").Value = 121.831933
$ws.Range("O5").Value = 0.9365419382944963
$ws.Range("P5").Value = 0.9365419382944962
$ws.Range("Q5").Value = 85.70199642477778
$ws.Range("R5").Value = 771.317967823
$ws.Range("S5").Value = 0.008242245193452144
$ws.Range("T5").Value = 0.008242245193452144
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pecam1"
$ws.Range("C6").Value = "Cd38"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.110333333333333
$ws.Range("H6").Value = 6.331
$ws.Range("I6").Value = 0.008800721950008783
$ws.Range("J6").Value = 0.008800721950008783
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.454549333333333
$ws.Range("N6").Value = 7.363648
$ws.Range("O6").Value = 0.05660556309845621
$ws.Range("P6").Value = 0.05660556309845621
$ws.Range("Q6").Value = 5.179917276444443
$ws.Range("R6").Value = 46.61925548799999
$ws.Range("S6").Value = 0.0004981698216531907
$ws.Range("T6").Value = 0.0004981698216531907
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pecam1"
$ws.Range("C7").Value = "Cd38"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.110333333333333
$ws.Range("H7").Value = 6.331
$ws.Range("I7").Value = 0.008800721950008783
$ws.Range("J7").Value = 0.008800721950008783
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2971403333333333
$ws.Range("N7").Value = 0.891421
$ws.Range("O7").Value = 0.006852498607047613
$ws.Range("P7").Value = 0.006852498607047612
$ws.Range("Q7").Value = 0.6270651501111111
$ws.Range("R7").Value = 5.643586351
$ws.Range("S7").Value = 0.00006030693490344854
$ws.Range("T7").Value = 0.00006030693490344853
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Pecam1"
$ws.Range("C8").Value = "Cd38"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.099218
$ws.Range("H8").Value = 9.297654
$ws.Range("I8").Value = 0.01292466713653245
$ws.Range("J8").Value = 0.01292466713653245
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 40.61064433333333
$ws.Range("N8").Value = 121.831933
$ws.Range("O8").Value = 0.9365419382944963
$ws.Range("P8").Value = 0.9365419382944962
$ws.Range("Q8").Value = 125.8612399094647
$ws.Range("R8").Value = 1132.751159185182
$ws.Range("S8").Value = 0.01210449281185928
$ws.Range("T8").Value = 0.01210449281185928
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Pecam1"
$ws.Range("C9").Value = "Cd38"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.099218
$ws.Range("H9").Value = 9.297654
$ws.Range("I9").Value = 0.01292466713653245
$ws.Range("J9").Value = 0.01292466713653245
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.454549333333333
$ws.Range("N9").Value = 7.363648
$ws.Range("O9").Value = 0.05660556309845621
$ws.Range("P9").Value = 0.05660556309845621
$ws.Range("Q9").Value = 7.607183475754666
$ws.Range("R9").Value = 68.46465128179199
$ws.Range("S9").Value = 0.0007316080611235312
$ws.Range("T9").Value = 0.0007316080611235311
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Pecam1"
$ws.Range("C10").Value = "Cd38"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.099218
$ws.Range("H10").Value = 9.297654
$ws.Range("I10").Value = 0.01292466713653245
$ws.Range("J10").Value = 0.01292466713653245
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.2971403333333333
$ws.Range("N10").Value = 0.891421
$ws.Range("O10").Value = 0.006852498607047613
$ws.Range("P10").Value = 0.006852498607047612
$ws.Range("Q10").Value = 0.9209026695926666
$ws.Range("R10").Value = 8.288124026334
$ws.Range("S10").Value = 0.0000885662635496427
$ws.Range("T10").Value = 0.00008856626354964269
